$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = "Начал взаимодействие с консультантом!"
$ws.Range("E7").Value = "11/11/2023 11:50:21"

# Row 8
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = 3
$ws.Range("D8").Value = "Успешно добавлен в базу!"
$ws.Range("E8").Value = "11/11/2023 16:08:54"

# Match the styling used by the other cells in column A (style index 1: bold font, thin box border, centered/top aligned)
$styleRange = $ws.Range("A7:A8")
$styleRange.Font.Bold = $true
$styleRange.HorizontalAlignment = -4108
$styleRange.VerticalAlignment = -4160
$styleRange.Borders.LineStyle = 1
$styleRange.Borders.Weight = 2
